# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps recorded on the Overview, zh-cn and
# de-de sheets to reflect a newer handback report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
# (this value is shared with de-de!H2, the "Correspond Handoff Datetime" for
# the same source file, and both must be kept in sync).
$wsOverview.Range("G2").Value = "2016-09-05 01:12:19"
$wsDeDe.Range("H2").Value     = "2016-09-05 01:12:19"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-05 01:12:14"
$wsZhCn.Range("K2").Value = "2016-09-05 01:12:31"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-09-05 01:12:39"
